$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Copy the formatting (border, alignment, number formats) from the last
# existing data row (69) down into the new row (70) before writing values,
# so the new row matches the established look of the table.
$ws.Range("A69:J69").Copy() | Out-Null
$ws.Range("A70:J70").PasteSpecial(-4122) | Out-Null

# Append the latest day's PTC run totals (date serial 42543 = 2016-06-22,
# reported "as of the 23rd" per the commit message).
$ws.Range("A70").Value = 42543
$ws.Range("B70").Value = 142
$ws.Range("C70").Value = 126
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 16
$ws.Range("F70").Formula = "=C70"
$ws.Range("G70").Formula = "=F70/B70"
$ws.Range("H70").Value = 44.351851851845694
$ws.Range("I70").Value = 35.966666660970077
$ws.Range("J70").Value = 56.366666662506759

$ws.Range("C70").Select() | Out-Null
